$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A506").Value = 505
$ws.Range("B506").Value = 'Saturday, Jan 14'
$ws.Range("C506").Value = '6:00 PM'
$ws.Range("D506").Value = 'FR6221'
$ws.Range("E506").Value = 'Bristol'
$ws.Range("F506").Value = '(BRS)'
$ws.Range("G506").Value = 'Buzz '
$ws.Range("H506").Value = 'B38M'
$ws.Range("I506").Value = '(SP-RZF)'
$ws.Range("J506").Value = '6:07 PM'
$ws.Range("L506").Value = '0 hours, 7 minutes'

$ws.Range("A507").Value = 506
$ws.Range("B507").Value = 'Saturday, Jan 14'
$ws.Range("C507").Value = '6:05 PM'
$ws.Range("D507").Value = 'BA872'
$ws.Range("E507").Value = 'London'
$ws.Range("F507").Value = '(LHR)'
$ws.Range("G507").Value = 'British Airways '
$ws.Range("H507").Value = 'A320'
$ws.Range("I507").Value = '(G-MIDO)'
$ws.Range("J507").Value = '6:05 PM'
$ws.Range("L507").Value = '0 hours, 0 minutes'

$ws.Range("A508").Value = 507
$ws.Range("B508").Value = 'Saturday, Jan 14'
$ws.Range("C508").Value = '6:30 PM'
$ws.Range("D508").Value = 'FR6876'
$ws.Range("E508").Value = 'Milan'
$ws.Range("F508").Value = '(BGY)'
$ws.Range("G508").Value = 'Buzz '
$ws.Range("H508").Value = 'B38M'
$ws.Range("I508").Value = '(SP-RZD)'
$ws.Range("J508").Value = '6:22 PM'
$ws.Range("L508").Value = '0 hours, -8 minutes'

$ws.Range("A509").Value = 508
$ws.Range("B509").Value = 'Saturday, Jan 14'
$ws.Range("C509").Value = '6:33 PM'
$ws.Range("D509").Value = 'RR9502'
$ws.Range("E509").Value = 'Tel Aviv'
$ws.Range("F509").Value = '(TLV)'
$ws.Range("G509").Value = 'Ryanair '
$ws.Range("H509").Value = 'B738'
$ws.Range("I509").Value = '(SP-RSH)'
$ws.Range("J509").Value = 'Diverted to KTW'

$ws.Range("A510").Value = 509
$ws.Range("B510").Value = 'Saturday, Jan 14'
$ws.Range("C510").Value = '6:35 PM'
$ws.Range("D510").Value = 'FR5891'
$ws.Range("E510").Value = 'Eindhoven'
$ws.Range("F510").Value = '(EIN)'
$ws.Range("G510").Value = 'Buzz '
$ws.Range("H510").Value = 'B38M'
$ws.Range("I510").Value = '(SP-RZC)'
$ws.Range("J510").Value = '6:52 PM'
$ws.Range("L510").Value = '0 hours, 17 minutes'

$ws.Range("A511").Value = 510
$ws.Range("B511").Value = 'Saturday, Jan 14'
$ws.Range("C511").Value = '6:40 PM'
$ws.Range("D511").Value = 'FR1642'
$ws.Range("E511").Value = 'Vienna'
$ws.Range("F511").Value = '(VIE)'
$ws.Range("G511").Value = 'Ryanair '
$ws.Range("H511").Value = 'B738'
$ws.Range("I511").Value = '(SP-RST)'
$ws.Range("J511").Value = '6:36 PM'
$ws.Range("L511").Value = '0 hours, -4 minutes'

$ws.Range("A512").Value = 511
$ws.Range("B512").Value = 'Saturday, Jan 14'
$ws.Range("C512").Value = '6:45 PM'
$ws.Range("D512").Value = 'FR3986'
$ws.Range("E512").Value = 'Bari'
$ws.Range("F512").Value = '(BRI)'
$ws.Range("G512").Value = 'Ryanair '
$ws.Range("H512").Value = 'B738'
$ws.Range("I512").Value = '(SP-RKB)'
$ws.Range("J512").Value = '6:42 PM'
$ws.Range("L512").Value = '0 hours, -3 minutes'

$ws.Range("A513").Value = 512
$ws.Range("B513").Value = 'Saturday, Jan 14'
$ws.Range("C513").Value = '7:25 PM'
$ws.Range("D513").Value = 'FR6277'
$ws.Range("E513").Value = 'Nuremberg'
$ws.Range("F513").Value = '(NUE)'
$ws.Range("G513").Value = 'Ryanair '
$ws.Range("H513").Value = 'B738'
$ws.Range("I513").Value = '(9H-QCX)'
$ws.Range("J513").Value = '6:56 PM'
$ws.Range("L513").Value = '0 hours, -29 minutes'

$ws.Range("A514").Value = 513
$ws.Range("B514").Value = 'Saturday, Jan 14'
$ws.Range("C514").Value = '7:32 PM'
$ws.Range("D514").Value = 'RR9504'
$ws.Range("E514").Value = 'Tel Aviv'
$ws.Range("F514").Value = '(TLV)'
$ws.Range("G514").Value = 'Ryanair '
$ws.Range("H514").Value = 'B738'
$ws.Range("I514").Value = '(SP-RSA)'
$ws.Range("J514").Value = '7:21 PM'
$ws.Range("L514").Value = '0 hours, -11 minutes'

$ws.Range("A515").Value = 514
$ws.Range("B515").Value = 'Saturday, Jan 14'
$ws.Range("C515").Value = '8:25 PM'
$ws.Range("D515").Value = 'FR7968'
$ws.Range("E515").Value = 'Pisa'
$ws.Range("F515").Value = '(PSA)'
$ws.Range("G515").Value = 'Ryanair '
$ws.Range("H515").Value = 'B738'
$ws.Range("I515").Value = '(9H-QDG)'
$ws.Range("J515").Value = '8:15 PM'
$ws.Range("L515").Value = '0 hours, -10 minutes'

$ws.Range("A516").Value = 515
$ws.Range("B516").Value = 'Saturday, Jan 14'
$ws.Range("C516").Value = '8:40 PM'
$ws.Range("D516").Value = 'FR6211'
$ws.Range("E516").Value = 'Paris'
$ws.Range("F516").Value = '(BVA)'
$ws.Range("G516").Value = 'Buzz '
$ws.Range("H516").Value = 'B38M'
$ws.Range("I516").Value = '(SP-RZB)'
$ws.Range("J516").Value = '8:43 PM'
$ws.Range("L516").Value = '0 hours, 3 minutes'

$ws.Range("A517").Value = 516
$ws.Range("B517").Value = 'Saturday, Jan 14'
$ws.Range("C517").Value = '8:40 PM'
$ws.Range("D517").Value = 'W65094'
$ws.Range("E517").Value = 'Abu Dhabi'
$ws.Range("F517").Value = '(AUH)'
$ws.Range("G517").Value = 'Wizz Air '
$ws.Range("H517").Value = 'A21N'
$ws.Range("I517").Value = '(HA-LVG)'
$ws.Range("J517").Value = '8:55 PM'
$ws.Range("L517").Value = '0 hours, 15 minutes'

$ws.Range("A518").Value = 517
$ws.Range("B518").Value = 'Saturday, Jan 14'
$ws.Range("C518").Value = '8:50 PM'
$ws.Range("D518").Value = 'FR2332'
$ws.Range("E518").Value = 'Leeds'
$ws.Range("F518").Value = '(LBA)'
$ws.Range("G518").Value = 'Ryanair '
$ws.Range("H518").Value = 'B738'
$ws.Range("I518").Value = '(EI-EKN)'
$ws.Range("J518").Value = '8:22 PM'
$ws.Range("L518").Value = '0 hours, -28 minutes'

$ws.Range("A519").Value = 518
$ws.Range("B519").Value = 'Saturday, Jan 14'
$ws.Range("C519").Value = '8:50 PM'
$ws.Range("D519").Value = 'FR2712'
$ws.Range("E519").Value = 'London'
$ws.Range("F519").Value = '(STN)'
$ws.Range("G519").Value = 'Ryanair '
$ws.Range("H519").Value = 'B738'
$ws.Range("I519").Value = '(EI-ENL)'
$ws.Range("J519").Value = '8:46 PM'
$ws.Range("L519").Value = '0 hours, -4 minutes'

$ws.Range("A520").Value = 519
$ws.Range("B520").Value = 'Saturday, Jan 14'
$ws.Range("C520").Value = '9:25 PM'
$ws.Range("D520").Value = 'FR6249'
$ws.Range("E520").Value = 'Manchester'
$ws.Range("F520").Value = '(MAN)'
$ws.Range("G520").Value = 'Ryanair '
$ws.Range("H520").Value = 'B738'
$ws.Range("I520").Value = '(EI-DHB)'
$ws.Range("J520").Value = '9:08 PM'
$ws.Range("L520").Value = '0 hours, -17 minutes'

$ws.Range("A521").Value = 520
$ws.Range("B521").Value = 'Saturday, Jan 14'
$ws.Range("C521").Value = '10:05 PM'
$ws.Range("D521").Value = 'FR2023'
$ws.Range("E521").Value = 'Dublin'
$ws.Range("F521").Value = '(DUB)'
$ws.Range("G521").Value = 'Buzz '
$ws.Range("H521").Value = 'B38M'
$ws.Range("I521").Value = '(SP-RZH)'
$ws.Range("J521").Value = '9:52 PM'
$ws.Range("L521").Value = '0 hours, -13 minutes'

$ws.Range("A522").Value = 521
$ws.Range("B522").Value = 'Saturday, Jan 14'
$ws.Range("C522").Value = '10:05 PM'
$ws.Range("D522").Value = 'LH1624'
$ws.Range("E522").Value = 'Munich'
$ws.Range("F522").Value = '(MUC)'
$ws.Range("G522").Value = 'Lufthansa '
$ws.Range("H522").Value = 'A320'
$ws.Range("I522").Value = '(D-AIWJ)'
$ws.Range("J522").Value = '10:09 PM'
$ws.Range("L522").Value = '0 hours, 4 minutes'

$ws.Range("A523").Value = 522
$ws.Range("B523").Value = 'Saturday, Jan 14'
$ws.Range("C523").Value = '10:10 PM'
$ws.Range("D523").Value = 'FR4934'
$ws.Range("E523").Value = 'Brussels'
$ws.Range("F523").Value = '(CRL)'
$ws.Range("G523").Value = 'Ryanair '
$ws.Range("H523").Value = 'B738'
$ws.Range("I523").Value = '(SP-RKC)'
$ws.Range("J523").Value = '10:06 PM'
$ws.Range("L523").Value = '0 hours, -4 minutes'

$ws.Range("A524").Value = 523
$ws.Range("B524").Value = 'Saturday, Jan 14'
$ws.Range("C524").Value = '10:40 PM'
$ws.Range("D524").Value = 'FR6243'
$ws.Range("E524").Value = 'Naples'
$ws.Range("F524").Value = '(NAP)'
$ws.Range("G524").Value = 'Buzz '
$ws.Range("H524").Value = 'B38M'
$ws.Range("I524").Value = '(SP-RZF)'
$ws.Range("J524").Value = '11:14 PM'
$ws.Range("L524").Value = '0 hours, 34 minutes'

$ws.Range("A525").Value = 524
$ws.Range("B525").Value = 'Saturday, Jan 14'
$ws.Range("C525").Value = '10:50 PM'
$ws.Range("D525").Value = 'KL1999'
$ws.Range("E525").Value = 'Amsterdam'
$ws.Range("F525").Value = '(AMS)'
$ws.Range("G525").Value = 'KLM '
$ws.Range("H525").Value = 'E295'
$ws.Range("I525").Value = '(PH-NXF)'
$ws.Range("J525").Value = '10:30 PM'
$ws.Range("L525").Value = '0 hours, -20 minutes'

$ws.Range("A526").Value = 525
$ws.Range("B526").Value = 'Saturday, Jan 14'
$ws.Range("C526").Value = '11:15 PM'
$ws.Range("D526").Value = 'W65082'
$ws.Range("E526").Value = 'Kutaisi'
$ws.Range("F526").Value = '(KUT)'
$ws.Range("G526").Value = 'Wizz Air '
$ws.Range("H526").Value = 'A21N'
$ws.Range("I526").Value = '(HA-LVO)'
$ws.Range("J526").Value = '10:55 PM'
$ws.Range("L526").Value = '0 hours, -20 minutes'

$ws.Range("A527").Value = 526
$ws.Range("B527").Value = 'Saturday, Jan 14'
$ws.Range("C527").Value = '11:20 PM'
$ws.Range("D527").Value = 'LH1370'
$ws.Range("E527").Value = 'Frankfurt'
$ws.Range("F527").Value = '(FRA)'
$ws.Range("G527").Value = 'Lufthansa '
$ws.Range("H527").Value = 'A321'
$ws.Range("I527").Value = '(D-AIDG)'
$ws.Range("J527").Value = '11:22 PM'
$ws.Range("L527").Value = '0 hours, 2 minutes'

$ws.Range("A528").Value = 527
$ws.Range("B528").Value = 'Saturday, Jan 14'
$ws.Range("C528").Value = '11:25 PM'
$ws.Range("D528").Value = 'FR1813'
$ws.Range("E528").Value = 'London'
$ws.Range("F528").Value = '(LTN)'
$ws.Range("G528").Value = 'Ryanair '
$ws.Range("H528").Value = 'B38M'
$ws.Range("I528").Value = '(SP-RZK)'
$ws.Range("J528").Value = '10:51 PM'
$ws.Range("L528").Value = '0 hours, -34 minutes'

$ws.Range("A529").Value = 528
$ws.Range("B529").Value = 'Saturday, Jan 14'
$ws.Range("C529").Value = '11:25 PM'
$ws.Range("D529").Value = 'FR6217'
$ws.Range("E529").Value = 'Oslo'
$ws.Range("F529").Value = '(TRF)'
$ws.Range("G529").Value = 'Buzz '
$ws.Range("H529").Value = 'B38M'
$ws.Range("I529").Value = '(SP-RZD)'
$ws.Range("J529").Value = '11:59 PM'
$ws.Range("L529").Value = '0 hours, 34 minutes'

$ws.Range("A530").Value = 529
$ws.Range("B530").Value = 'Saturday, Jan 14'
$ws.Range("C530").Value = '11:25 PM'
$ws.Range("D530").Value = 'W65004'
$ws.Range("E530").Value = 'London'
$ws.Range("F530").Value = '(LTN)'
$ws.Range("G530").Value = 'Wizz Air '
$ws.Range("H530").Value = 'A321'
$ws.Range("I530").Value = '(HA-LXO)'
$ws.Range("J530").Value = '11:11 PM'
$ws.Range("L530").Value = '0 hours, -14 minutes'

$ws.Range("A531").Value = 530
$ws.Range("B531").Value = 'Saturday, Jan 14'
$ws.Range("C531").Value = '11:35 PM'
$ws.Range("D531").Value = 'FR6245'
$ws.Range("E531").Value = 'Lille'
$ws.Range("F531").Value = '(LIL)'
$ws.Range("G531").Value = 'Buzz '
$ws.Range("H531").Value = 'B38M'
$ws.Range("I531").Value = '(SP-RZC)'
$ws.Range("J531").Value = '11:49 PM'
$ws.Range("L531").Value = '0 hours, 14 minutes'

$ws.Range("A532").Value = 531
$ws.Range("B532").Value = 'Saturday, Jan 14'
$ws.Range("C532").Value = '11:35 PM'
$ws.Range("D532").Value = 'FR6257'
$ws.Range("E532").Value = 'Stockholm'
$ws.Range("F532").Value = '(ARN)'
$ws.Range("G532").Value = 'Ryanair '
$ws.Range("H532").Value = 'B738'
$ws.Range("I532").Value = '(SP-RKB)'
$ws.Range("J532").Value = '11:44 PM'
$ws.Range("L532").Value = '0 hours, 9 minutes'

$ws.Range("A533").Value = 532
$ws.Range("B533").Value = 'Saturday, Jan 14'
$ws.Range("C533").Value = '11:35 PM'
$ws.Range("D533").Value = 'LO3911'
$ws.Range("E533").Value = 'Warsaw'
$ws.Range("F533").Value = '(WAW)'
$ws.Range("G533").Value = 'LOT '
$ws.Range("H533").Value = 'E190'
$ws.Range("I533").Value = '(SP-LMH)'
$ws.Range("J533").Value = '11:17 PM'
$ws.Range("L533").Value = '0 hours, -18 minutes'

$ws.Range("A534").Value = 533
$ws.Range("B534").Value = 'Saturday, Jan 14'
$ws.Range("C534").Value = '11:50 PM'
$ws.Range("D534").Value = 'FR6231'
$ws.Range("E534").Value = 'Szczecin'
$ws.Range("F534").Value = '(SZZ)'
$ws.Range("G534").Value = 'Buzz '
$ws.Range("H534").Value = 'B38M'
$ws.Range("I534").Value = '(SP-RZB)'
$ws.Range("J534").Value = '12:25 AM'
$ws.Range("L534").Value = '0 hours, 35 minutes'

$ws.Range("A535").Value = 534
$ws.Range("B535").Value = 'Saturday, Jan 14'
$ws.Range("C535").Value = '11:50 PM'
$ws.Range("D535").Value = 'W65062'
$ws.Range("E535").Value = 'Rome'
$ws.Range("F535").Value = '(FCO)'
$ws.Range("G535").Value = 'Wizz Air '
$ws.Range("H535").Value = 'A21N'
$ws.Range("I535").Value = '(HA-LZI)'
$ws.Range("J535").Value = '11:41 PM'
$ws.Range("L535").Value = '0 hours, -9 minutes'

$ws.Range("K506:K535").Font.Size = 11
$ws.Range("M506:M535").Font.Size = 11
$ws.Range("L509").Font.Size = 11
